$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.063.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.290.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.651"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.36%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.125"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.406"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.865.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "66.114.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000163"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.280.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "433.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.433.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.510"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000113"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.196"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("E38").Value = "  -5.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.781.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.771"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0659"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "316.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0268"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.94%  "
